$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A/B/C columns one slot right -> B/C/D. This carries the
# pre-existing column-C width formatting (width=12, bestFit, customWidth)
# along to column D automatically, matching the author's edit.
$ws.Columns.Item(1).Insert()

# The "Actual" DB value used to live in A2; after the shift it is in B2.
# Rebind the BridgeDBValue defined name to that new location. (Deleting and
# re-creating the name - rather than editing RefersTo in place - is required
# for downstream formulas to actually pick up the new target cell.)
$wb.Names.Item(1).Delete()
$ws.Range("B2").Name = "BridgeDBValue"

# Row labels (set Actual/Reported before the new header so the shared-string
# table order matches the source order of insertion: Feet, Inches, Actual,
# Reported, From DB).
$ws.Range("A2").Value = "Actual"
$ws.Range("A3").Value = "Reported"
$ws.Range("B1").Value = "From DB"
$ws.Range("A5").Value = "Actual"
$ws.Range("A6").Value = "Reported"

# Row 3: "Reported" = the row-2 (Actual) conversion minus 3 inches.
$ws.Range("C3").Formula = "=(ROUND(BridgeDBValue/100,0)*12+MOD(BridgeDBValue,100))-3"
$ws.Range("D3").Formula = "=(ROUND(BridgeDBValue/100,0)+(MOD(BridgeDBValue,100)/12))-3/12"

# Row 5: second "Actual" DB record (its own literal value in B5, not the
# named range).
$ws.Range("B5").Value = 1400
$ws.Range("C5").Formula = "=ROUND(B5/100,0)*12+MOD(BridgeDBValue,100)"
$ws.Range("D5").Formula = "=ROUND(B5/100,0)+(MOD(B5,100)/12)"

# Row 6: "Reported" = row-5 conversion minus 3 inches.
$ws.Range("C6").Formula = "=(ROUND(B5/100,0)*12+MOD(B5,100))-3"
$ws.Range("D6").Formula = "=(ROUND(B5/100,0)+(MOD(B5,100)/12))-3/12"

# Apply the built-in "Calculation" cell style (orange bold font, light-grey
# fill, thin grey border) to every computed conversion cell.
$ws.Range("C2:D3").Style = "Calculation"
$ws.Range("C5:D6").Style = "Calculation"

# Match the author's final selection.
$ws.Range("C5").Select() | Out-Null
